$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (row 2 through 67) date serial value from 45184 to 45186
for ($r = 2; $r -le 67; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value2 = 45186
    }
}

# Update HYPERLINK formulas in row 2 (columns S, T, V, W, X, Y) to include friendly name text
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_GRASTORP/artfynd/A 35302-2021.xlsx", "A 35302-2021")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_GRASTORP/kartor/A 35302-2021.png", "A 35302-2021")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_GRASTORP/klagomål/A 35302-2021.docx", "A 35302-2021")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_GRASTORP/klagomålsmail/A 35302-2021.docx", "A 35302-2021")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_GRASTORP/tillsyn/A 35302-2021.docx", "A 35302-2021")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_GRASTORP/tillsynsmail/A 35302-2021.docx", "A 35302-2021")'
